# [TEST SCRAPE] updated files from azure vm
#
# 1) "ODI Batting" sheet: remove the stray empty cell B2.
# 2) Add a new worksheet "ODI Batting Extra" (after "ODI Bowling") with a
#    header row (MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6,
#    PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH) styled like the other sheets'
#    header rows, and a single data row seeded with MATCH_CODE 4602.

$wb = $excel.ActiveWorkbook

# --- 1. Remove the empty B2 cell on "ODI Batting" ---------------------
$battingWs = $wb.Worksheets.Item("ODI Batting")
$battingWs.Range("B2").ClearContents()

# --- 2. Add the new "ODI Batting Extra" sheet after the last sheet ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Header row values
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Match the bold/centered/bordered header styling used on the other sheets
# by copying the format from an existing header cell.
$styleSource = $wb.Worksheets.Item("Player Info").Range("A1")
$styleSource.Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)

# Data row: MATCH_CODE = 4602, remaining columns present but blank.
# Force text typing (so "4602" stays a string, not a number) via a
# temporary text number format, then restore the default style so no
# extra formatting is left behind on the data row.
$dataRange = $newSheet.Range("A2:F2")
$dataRange.NumberFormat = "@"
$newSheet.Cells.Item(2, 1).Value = "4602"
$newSheet.Cells.Item(2, 2).Value = ""
$newSheet.Cells.Item(2, 3).Value = ""
$newSheet.Cells.Item(2, 4).Value = ""
$newSheet.Cells.Item(2, 5).Value = ""
$newSheet.Cells.Item(2, 6).Value = ""
$dataRange.Style = "Normal"

# Restore original active sheet/selection.
$wb.Worksheets.Item("Player Info").Activate()

Write-Output "Applied ODI Batting Extra sheet + ODI Batting!B2 cleanup"
